$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3854.8
$ws.Range("I64").Value = 3972
$ws.Range("J64").Value = 3776.6667
$ws.Range("K64").Value = 3972
$ws.Range("L64").Value = 3776.6667
$ws.Range("M64").Value = -3724
$ws.Range("N64").Value = -4272.6667
$ws.Range("H67").Value = 3854.8
$ws.Range("I67").Value = 3972
$ws.Range("J67").Value = 3776.6667
$ws.Range("K67").Value = 3972
$ws.Range("L67").Value = 3776.6667
$ws.Range("M67").Value = -3114
$ws.Range("N67").Value = -5492.6667
$ws.Range("H69").Value = 3337.6667
$ws.Range("I69").Value = 2013
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 6039
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -5165
$ws.Range("N69").Value = -13748
$ws.Range("H72").Value = 3337.6667
$ws.Range("I72").Value = 2013
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 18117
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -13749
$ws.Range("N72").Value = -44736
$ws.Range("H124").Value = 39800
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 39800
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 39800
$ws.Range("N124").Value = -49620

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 100001544
$ws.Range("I61").Value = 125001304
$ws.Range("J61").Value = 2507
$ws.Range("K61").Value = 125001304
$ws.Range("L61").Value = 2507
$ws.Range("M61").Value = -125001092
$ws.Range("N61").Value = -2931
$ws.Range("H74").Value = 2753.5715
$ws.Range("I74").Value = 2143.9092
$ws.Range("J74").Value = 3424.2
$ws.Range("K74").Value = 2143.9092
$ws.Range("L74").Value = 3424.2
$ws.Range("M74").Value = -1269.9092
$ws.Range("N74").Value = -5172.2
$ws.Range("H77").Value = 2753.5715
$ws.Range("I77").Value = 2143.9092
$ws.Range("J77").Value = 3424.2
$ws.Range("K77").Value = 10719.546
$ws.Range("L77").Value = 17121
$ws.Range("M77").Value = -6351.546
$ws.Range("N77").Value = -25857
$ws.Range("H122").Value = 2248.6316
$ws.Range("I122").Value = 1965.2307
$ws.Range("J122").Value = 2862.6667
$ws.Range("K122").Value = 5895.6921
$ws.Range("L122").Value = 8588.000100000001
$ws.Range("M122").Value = -3445.6921
$ws.Range("N122").Value = -13488.0001
$ws.Range("H132").Value = 2783.7297
$ws.Range("I132").Value = 2377.45
$ws.Range("J132").Value = 3261.7058
$ws.Range("K132").Value = 7132.349999999999
$ws.Range("L132").Value = 9785.117400000001
$ws.Range("M132").Value = -4602.349999999999
$ws.Range("N132").Value = -14845.1174
$ws.Range("H136").Value = 100001544
$ws.Range("I136").Value = 125001304
$ws.Range("J136").Value = 2507
$ws.Range("K136").Value = 375003912
$ws.Range("L136").Value = 7521
$ws.Range("M136").Value = -375001362
$ws.Range("N136").Value = -12621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 3600
$ws.Range("J14").Value = 3600
$ws.Range("L14").Value = 3600
$ws.Range("N14").Value = -3940
$ws.Range("H15").Value = 1004.5
$ws.Range("J15").Value = 1004.5
$ws.Range("L15").Value = 1004.5
$ws.Range("N15").Value = -1344.5
$ws.Range("H21").Value = 3500
$ws.Range("J21").Value = 3500
$ws.Range("L21").Value = 3500
$ws.Range("N21").Value = -3970
$ws.Range("H23").Value = 4200
$ws.Range("J23").Value = 4200
$ws.Range("L23").Value = 4200
$ws.Range("N23").Value = -4680
$ws.Range("H26").Value = 5750
$ws.Range("J26").Value = 5750
$ws.Range("L26").Value = 5750
$ws.Range("N26").Value = -6324
$ws.Range("H27").Value = 4200
$ws.Range("J27").Value = 4200
$ws.Range("L27").Value = 4200
$ws.Range("N27").Value = -4584
$ws.Range("H31").Value = 1277.5532
$ws.Range("I31").Value = 1277.5532
$ws.Range("K31").Value = 1277.5532
$ws.Range("M31").Value = -982.5532000000001
$ws.Range("H33").Value = 600
$ws.Range("I33").Value = 600
$ws.Range("K33").Value = 600
$ws.Range("M33").Value = -221
$ws.Range("H34").Value = 1277.5532
$ws.Range("I34").Value = 1277.5532
$ws.Range("K34").Value = 1277.5532
$ws.Range("M34").Value = -1075.5532
$ws.Range("H35").Value = 262.5
$ws.Range("I35").Value = 262.5
$ws.Range("K35").Value = 262.5
$ws.Range("M35").Value = 31.5
$ws.Range("H41").Value = 25000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H42").Value = 4000
$ws.Range("I42").Value = 4000
$ws.Range("K42").Value = 4000
$ws.Range("M42").Value = -3407

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 83334400
$ws.Range("I129").Value = 333333340
$ws.Range("J129").Value = 20834664
$ws.Range("K129").Value = 1000000020
$ws.Range("L129").Value = 62503992
$ws.Range("M129").Value = -999995020
$ws.Range("N129").Value = -62513992
$ws.Range("H131").Value = 28575114
$ws.Range("J131").Value = 5153.5415
$ws.Range("L131").Value = 15460.6245
$ws.Range("N131").Value = -25540.6245

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4118.7144
$ws.Range("I122").Value = 4093.375
$ws.Range("J122").Value = 4199.8
$ws.Range("K122").Value = 12280.125
$ws.Range("L122").Value = 12599.4
$ws.Range("M122").Value = -9830.125
$ws.Range("N122").Value = -17499.4
$ws.Range("H132").Value = 2462.634
$ws.Range("I132").Value = 2453.25
$ws.Range("J132").Value = 2475.8823
$ws.Range("K132").Value = 7359.75
$ws.Range("L132").Value = 7427.646900000001
$ws.Range("M132").Value = -4829.75
$ws.Range("N132").Value = -12487.6469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 62501628
$ws.Range("I122").Value = 250000000
$ws.Range("J122").Value = 2168.3333
$ws.Range("K122").Value = 750000000
$ws.Range("L122").Value = 6504.999899999999
$ws.Range("M122").Value = -749997550
$ws.Range("N122").Value = -11404.9999
$ws.Range("H132").Value = 2927.25
$ws.Range("I132").Value = 2692.077
$ws.Range("J132").Value = 3364
$ws.Range("K132").Value = 8076.231000000001
$ws.Range("L132").Value = 10092
$ws.Range("M132").Value = -5546.231000000001
$ws.Range("N132").Value = -15152

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 500
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

Write-Output "edits applied"